$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("Z100").Value = "In Translation"
$ws.Range("Z100").Value = ""
